$d = $word.ActiveDocument

# Locate the paragraph "Going off one side (left or right) of the screen
# makes you appear on the other side." -- the new content is inserted right
# after it, ahead of the two trailing blank paragraphs that close out the
# document.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Going off one side*other side.*") {
        $targetIndex = $i
        break
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Going off one side' paragraph"
}

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $null = $paragraph.Range.InsertXML($xml)
}

$target = $d.Paragraphs($targetIndex)

# 1) A blank paragraph right after the target paragraph.
$null = $target.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs($targetIndex + 1)
Set-ParagraphXml $blankPara "<w:p/>"

# 2) The "Up arrow objects..." paragraph.
$null = $blankPara.Range.InsertParagraphAfter()
$upArrowPara = $d.Paragraphs($targetIndex + 2)
Set-ParagraphXml $upArrowPara @'
<w:p>
<w:r>
<w:t>Up arrow objects in the game are placed on the level sections.  Touching them triggers the level to scroll up to the next section.  Or it makes a button appear on the phones of the player&#8217;s team that touched it that they can press when they are ready to scroll.</w:t>
</w:r>
</w:p>
'@

# 3) The "Could have some down arrows too..." paragraph, with the
#    gramStart/gramEnd proofing-error markers around the sentence and its
#    trailing period split into its own run (as produced by Word's grammar
#    checker).
$null = $upArrowPara.Range.InsertParagraphAfter()
$downArrowPara = $d.Paragraphs($targetIndex + 3)
Set-ParagraphXml $downArrowPara @'
<w:p>
<w:proofErr w:type="gramStart"/>
<w:r>
<w:t>Could have some down arrows too, that would knock people off the top of the screen</w:t>
</w:r>
<w:r>
<w:t>.</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
</w:p>
'@
